# Applies the "Updated cryptos list" GitHub Actions refresh to Sheet1.
# For each changed cell we set the new literal text. Price-column (D) values
# that look numeric are written with a leading apostrophe so Excel keeps them
# as plain text (matching the original inlineStr cells), and the quote-prefix
# style that the apostrophe triggers is immediately reset back to Normal so no
# stray cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '51.112.29'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.959.41'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +0.65%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.16%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '380.28'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +1.27%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '102.17'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +0.50%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +1.52%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.592'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +1.23%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '36.39'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +0.34%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -1.30%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.0858'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +1.36%  '; ForceText = $false }
    @{ Cell = 'B13'; Value = 'Polkadot'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; ForceText = $false }
    @{ Cell = 'D13'; Value = '7.83'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +5.62%  '; ForceText = $false }
    @{ Cell = 'B14'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D14'; Value = '3.424.01'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  +0.44%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '18.31'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +2.00%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '2.967.24'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +0.90%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '11.27'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +3.30%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '0.996'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  +1.57%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '51.183.63'; ForceText = $false }
    @{ Cell = 'E19'; Value = '  +0.39%  '; ForceText = $false }
    @{ Cell = 'E20'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '12.36'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -1.55%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '0.0₃0962'; ForceText = $false }
    @{ Cell = 'E22'; Value = '  +0.58%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '70.42'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +2.74%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '3.27'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +2.78%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '266.68'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +1.04%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '7.83'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -3.08%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '7.20'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -7.02%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '25.85'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +0.96%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.165'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.70%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '0.110'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +3.69%  '; ForceText = $false }
    @{ Cell = 'B33'; Value = 'InjectiveProtocol'; ForceText = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false }
    @{ Cell = 'D33'; Value = '34.36'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +3.41%  '; ForceText = $false }
    @{ Cell = 'B34'; Value = 'OKB'; ForceText = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; ForceText = $false }
    @{ Cell = 'D34'; Value = '51.18'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +1.06%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '2.05'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +1.61%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.0435'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -1.48%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '3.22'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +4.00%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.117'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +0.68%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  +2.50%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '16.46'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +1.31%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '125.22'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +3.68%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  -0.54%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '3.54'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +7.66%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '21.44'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.25%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  +0.34%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '0.272'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -2.39%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.37'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +2.59%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.038.70'; ForceText = $false }
    @{ Cell = 'E49'; Value = '  +1.88%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.0320'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -4.56%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +5.57%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.ForceText) {
        $value = "'" + $value
    }
    $range = $ws.Range($u.Cell)
    $range.Value = $value
    if ($u.ForceText) {
        $range.Style = 'Normal'
    }
}
